$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$eLab2 = @'
1: The peering between as2r2 and 20.30.1.2 is not up.
2: named on local is running but answered with REFUSED when quering for .
3: `resolv.conf` file not found for device `as1r1`
4: `resolv.conf` file not found for device `as1r2`
5: `resolv.conf` file not found for device `as2r1`
6: `resolv.conf` file not found for device `as2r2`
7: `resolv.conf` file not found for device `as3r1`
8: `resolv.conf` file not found for device `as3r2`
9: ping: pc.net: Temporary failure in name resolution
10: ping: pc.net: Temporary failure in name resolution
11: ping: pc.net: Temporary failure in name resolution
12: ping: pc.net: Temporary failure in name resolution
13: ping: pc.net: Temporary failure in name resolution
14: ping: pc.net: Temporary failure in name resolution

'@

$eLab3 = @'
1: Devices connected to collision domain A ['as1r1', 'pc'] are different from the one in the template ['root', 'as1r2'].
2: Devices connected to collision domain K ['as1r2', 'local'] are different from the one in the template ['as3r2', 'local'].
3: Devices connected to collision domain J ['as3r2', 'root'] are different from the one in the template ['as3r2', 'pc'].
4: as1r2.startup file not found
5: as2r1.startup file not found
6: The interface `eth0` of `as1r1` has the following IP addresses: ['1.0.0.1/24']`.
7: The interface `eth1` of `as1r1` has the following IP addresses: ['10.20.0.1/30']`.
8: The interface `eth0` of `as1r2` has the following IP addresses: []`.
9: The interface `eth1` of `as1r2` has the following IP addresses: []`.
10: The interface `eth2` of `as1r2` has the following IP addresses: []`.
11: The interface `eth0` of `as2r1` has the following IP addresses: []`.
12: The interface `eth1` of `as2r1` has the following IP addresses: []`.
13: The interface `eth2` of `as2r1` has the following IP addresses: []`.
14: The interface `eth1` of `as2r2` has the following IP addresses: ['2.0.0.1/24']`.
15: The interface `eth2` of `as3r2` has the following IP addresses: ['3.1.0.1/24']`.
16: Interface eth`3` not found on `as3r2`
17: The interface `eth0` of `root` has the following IP addresses: ['3.1.0.2/24']`.
18: The interface `eth0` of `local` has the following IP addresses: ['1.2.0.2/24']`.
19: The interface `eth0` of `pc` has the following IP addresses: ['1.0.0.2/24']`.
20: No answer from `1.1.0.2` to `as1r1`.
21: ping: connect: Network is unreachable
22: ping: connect: Network is unreachable
23: ping: connect: Network is unreachable
24: ping: connect: Network is unreachable
25: ping: connect: Network is unreachable
26: ping: connect: Network is unreachable
27: ping: connect: Network is unreachable
28: ping: connect: Network is unreachable
29: ping: connect: Network is unreachable
30: ping: connect: Network is unreachable
31: No answer from `10.20.0.2` to `as1r1`.
32: ping: connect: Network is unreachable
33: ping: connect: Network is unreachable
34: ping: connect: Network is unreachable
35: ping: connect: Network is unreachable
36: ping: connect: Network is unreachable
37: ping: connect: Network is unreachable
38: ping: connect: Network is unreachable
39: ping: connect: Network is unreachable
40: ping: connect: Network is unreachable
41: ping: connect: Network is unreachable
42: ping: connect: Network is unreachable
43: ping: connect: Network is unreachable
44: ping: connect: Network is unreachable
45: ping: connect: Network is unreachable
46: ping: connect: Network is unreachable
47: ping: connect: Network is unreachable
48: ping: connect: Network is unreachable
49: ping: connect: Network is unreachable
50: ping: connect: Network is unreachable
51: ping: connect: Network is unreachable
52: ping: connect: Network is unreachable
53: ping: connect: Network is unreachable
54: ping: connect: Network is unreachable
55: ping: connect: Network is unreachable
56: ping: connect: Network is unreachable
57: ping: connect: Network is unreachable
58: ping: connect: Network is unreachable
59: ping: connect: Network is unreachable
60: ping: connect: Network is unreachable
61: ping: connect: Network is unreachable
62: ping: connect: Network is unreachable
63: ping: connect: Network is unreachable
64: ping: connect: Network is unreachable
65: ping: connect: Network is unreachable
66: ping: connect: Network is unreachable
67: ping: connect: Network is unreachable
68: ping: connect: Network is unreachable
69: ping: connect: Network is unreachable
70: ping: connect: Network is unreachable
71: ping: connect: Network is unreachable
72: ping: connect: Network is unreachable
73: ping: connect: Network is unreachable
74: ping: connect: Network is unreachable
75: ping: connect: Network is unreachable
76: ping: connect: Network is unreachable
77: ping: connect: Network is unreachable
78: ping: connect: Network is unreachable
79: ping: connect: Network is unreachable
80: ping: connect: Network is unreachable
81: ping: connect: Network is unreachable
82: ping: connect: Network is unreachable
83: ping: connect: Network is unreachable
84: ping: connect: Network is unreachable
85: ping: connect: Network is unreachable
86: No answer from `2.0.0.2` to `as2r2`.
87: ping: connect: Network is unreachable
88: ping: connect: Network is unreachable
89: ping: connect: Network is unreachable
90: ping: connect: Network is unreachable
91: ping: connect: Network is unreachable
92: ping: connect: Network is unreachable
93: ping: connect: Network is unreachable
94: ping: connect: Network is unreachable
95: No answer from `10.20.1.1` to `as2r2`.
96: ping: connect: Network is unreachable
97: ping: connect: Network is unreachable
98: ping: connect: Network is unreachable
99: ping: connect: Network is unreachable
100: ping: connect: Network is unreachable
101: ping: connect: Network is unreachable
102: ping: connect: Network is unreachable
103: ping: connect: Network is unreachable
104: ping: connect: Network is unreachable
105: ping: connect: Network is unreachable
106: ping: connect: Network is unreachable
107: ping: connect: Network is unreachable
108: ping: connect: Network is unreachable
109: ping: connect: Network is unreachable
110: ping: connect: Network is unreachable
111: ping: connect: Network is unreachable
112: ping: connect: Network is unreachable
113: ping: connect: Network is unreachable
114: No answer from `20.30.0.1` to `as3r1`.
115: ping: connect: Network is unreachable
116: ping: connect: Network is unreachable
117: ping: connect: Network is unreachable
118: ping: connect: Network is unreachable
119: ping: connect: Network is unreachable
120: ping: connect: Network is unreachable
121: ping: connect: Network is unreachable
122: ping: connect: Network is unreachable
123: ping: connect: Network is unreachable
124: ping: connect: Network is unreachable
125: ping: connect: Network is unreachable
126: ping: connect: Network is unreachable
127: ping: connect: Network is unreachable
128: ping: connect: Network is unreachable
129: ping: connect: Network is unreachable
130: ping: connect: Network is unreachable
131: ping: connect: Network is unreachable
132: ping: connect: Network is unreachable
133: ping: connect: Network is unreachable
134: ping: connect: Network is unreachable
135: ping: connect: Network is unreachable
136: ping: connect: Network is unreachable
137: ping: connect: Network is unreachable
138: ping: connect: Network is unreachable
139: ping: connect: Network is unreachable
140: ping: connect: Network is unreachable
141: ping: connect: Network is unreachable
142: ping: connect: Network is unreachable
143: ping: connect: Network is unreachable
144: ping: connect: Network is unreachable
145: ping: connect: Network is unreachable
146: ping: connect: Network is unreachable
147: ping: connect: Network is unreachable
148: ping: connect: Network is unreachable
149: ping: connect: Network is unreachable
150: ping: connect: Network is unreachable
151: ping: connect: Network is unreachable
152: ping: connect: Network is unreachable
153: ping: connect: Network is unreachable
154: ping: connect: Network is unreachable
155: ping: connect: Network is unreachable
156: ping: connect: Network is unreachable
157: ping: connect: Network is unreachable
158: ping: connect: Network is unreachable
159: ping: connect: Network is unreachable
160: ping: connect: Network is unreachable
161: ping: connect: Network is unreachable
162: ping: connect: Network is unreachable
163: ping: connect: Network is unreachable
164: ping: connect: Network is unreachable
165: ping: connect: Network is unreachable
166: ping: connect: Network is unreachable
167: ping: connect: Network is unreachable
168: ping: connect: Network is unreachable
169: ping: connect: Network is unreachable
170: ping: connect: Network is unreachable
171: ping: connect: Network is unreachable
172: ping: connect: Network is unreachable
173: ping: connect: Network is unreachable
174: ping: connect: Network is unreachable
175: ping: connect: Network is unreachable
176: ping: connect: Network is unreachable
177: ping: connect: Network is unreachable
178: ping: connect: Network is unreachable
179: ping: connect: Network is unreachable
180: ping: connect: Network is unreachable
181: ping: connect: Network is unreachable
182: ping: connect: Network is unreachable
183: ping: connect: Network is unreachable
184: ping: connect: Network is unreachable
185: ping: connect: Network is unreachable
186: ping: connect: Network is unreachable
187: ping: connect: Network is unreachable
188: ping: connect: Network is unreachable
189: ping: connect: Network is unreachable
190: ping: connect: Network is unreachable
191: ping: connect: Network is unreachable
192: ping: connect: Network is unreachable
193: ping: connect: Network is unreachable
194: ping: connect: Network is unreachable
195: ping: connect: Network is unreachable
196: ping: connect: Network is unreachable
197: ping: connect: Network is unreachable
198: ping: connect: Network is unreachable
199: ping: connect: Network is unreachable
200: ping: connect: Network is unreachable
201: ping: connect: Network is unreachable
202: ping: connect: Network is unreachable
203: ping: connect: Network is unreachable
204: ping: connect: Network is unreachable
205: ping: connect: Network is unreachable
206: ping: connect: Network is unreachable
207: ping: connect: Network is unreachable
208: ping: connect: Network is unreachable
209: ping: connect: Network is unreachable
210: ping: connect: Network is unreachable
211: ping: connect: Network is unreachable
212: ping: connect: Network is unreachable
213: ping: connect: Network is unreachable
214: ping: connect: Network is unreachable
215: Daemon bgpd is not running on device `as1r1`
216: Daemon ripd is not running on device `as1r1`
217: Daemon bgpd is not running on device `as1r2`
218: Daemon ripd is not running on device `as1r2`
219: Daemon bgpd is not running on device `as2r1`
220: Daemon ripd is not running on device `as2r1`
221: Daemon bgpd is not running on device `as2r2`
222: Daemon ripd is not running on device `as2r2`
223: Daemon bgpd is not running on device `as3r1`
224: Daemon ripd is not running on device `as3r1`
225: Daemon bgpd is not running on device `as3r2`
226: Daemon ripd is not running on device `as3r2`
227: Daemon named is not running on device `local`
228: Daemon named is not running on device `net`
229: Daemon named is not running on device `root`
230: ERROR: bgpd is not running

231: ERROR: bgpd is not running

232: ERROR: Exiting: failed to connect to any daemons.

233: ERROR: Exiting: failed to connect to any daemons.

234: ERROR: Exiting: failed to connect to any daemons.

235: ERROR: Exiting: failed to connect to any daemons.

236: ERROR: Exiting: failed to connect to any daemons.

237: ERROR: bgpd is not running

238: ERROR: bgpd is not running

239: ERROR: bgpd is not running

240: ERROR: bgpd is not running

241: ERROR: bgpd is not running

242: ERROR: bgpd is not running

243: ERROR: bgpd is not running

244: Network 1.0.0.0/8 is not announced in bgpd.
245: Network 1.0.0.0/8 is not announced in bgpd.
246: Network 2.0.0.0/8 is not announced in bgpd.
247: Network 2.0.0.0/8 is not announced in bgpd.
248: Network 3.0.0.0/8 is not announced in bgpd.
249: Network 3.0.0.0/8 is not announced in bgpd.
250: connected routes are not injected into `ripd` on `as1r1`.
251: bgp routes are not injected into `ripd` on `as1r1`.
252: connected routes are not injected into `ripd` on `as1r2`.
253: bgp routes are not injected into `ripd` on `as1r2`.
254: connected routes are not injected into `ripd` on `as2r1`.
255: bgp routes are not injected into `ripd` on `as2r1`.
256: connected routes are not injected into `ripd` on `as2r2`.
257: bgp routes are not injected into `ripd` on `as2r2`.
258: connected routes are not injected into `ripd` on `as3r1`.
259: bgp routes are not injected into `ripd` on `as3r1`.
260: connected routes are not injected into `ripd` on `as3r2`.
261: bgp routes are not injected into `ripd` on `as3r2`.
262: The route 2.0.0.0/8 IS NOT found in the routing table of `as1r1`.
263: The route 3.0.0.0/8 IS NOT found in the routing table of `as1r1`.
264: The route 10.20.1.0/30 IS NOT found in the routing table of `as1r1`.
265: The route 20.30.0.0/30 IS NOT found in the routing table of `as1r1`.
266: The route 20.30.1.0/30 IS NOT found in the routing table of `as1r1`.
267: The route 1.0.0.0/24 IS NOT found in the routing table of `as1r2`.
268: The route 1.1.0.0/24 IS NOT found in the routing table of `as1r2`.
269: The route 2.0.0.0/8 IS NOT found in the routing table of `as1r2`.
270: The route 3.0.0.0/8 IS NOT found in the routing table of `as1r2`.
271: The route 10.20.0.0/30 IS NOT found in the routing table of `as1r2`.
272: The route 10.20.1.0/30 IS NOT found in the routing table of `as1r2`.
273: The route 20.30.0.0/30 IS NOT found in the routing table of `as1r2`.
274: The route 20.30.1.0/30 IS NOT found in the routing table of `as1r2`.
275: The route 1.0.0.0/8 IS NOT found in the routing table of `as2r1`.
276: The route 2.0.0.0/24 IS NOT found in the routing table of `as2r1`.
277: The route 2.1.0.0/24 IS NOT found in the routing table of `as2r1`.
278: The route 3.0.0.0/8 IS NOT found in the routing table of `as2r1`.
279: The route 10.20.0.0/30 IS NOT found in the routing table of `as2r1`.
280: The route 10.20.1.0/30 IS NOT found in the routing table of `as2r1`.
281: The route 20.30.0.0/30 IS NOT found in the routing table of `as2r1`.
282: The route 20.30.1.0/30 IS NOT found in the routing table of `as2r1`.
283: The route 1.0.0.0/8 IS NOT found in the routing table of `as2r2`.
284: The route 3.0.0.0/8 IS NOT found in the routing table of `as2r2`.
285: The route 10.20.0.0/30 IS NOT found in the routing table of `as2r2`.
286: The route 20.30.0.0/30 IS NOT found in the routing table of `as2r2`.
287: The route 1.0.0.0/8 IS NOT found in the routing table of `as3r1`.
288: The route 2.0.0.0/8 IS NOT found in the routing table of `as3r1`.
289: The route 3.1.0.0/24 IS NOT found in the routing table of `as3r1`.
290: The route 3.2.0.0/24 IS NOT found in the routing table of `as3r1`.
291: The route 10.20.0.0/30 IS NOT found in the routing table of `as3r1`.
292: The route 10.20.1.0/30 IS NOT found in the routing table of `as3r1`.
293: The route 20.30.1.0/30 IS NOT found in the routing table of `as3r1`.
294: The route 0.0.0.0/0 IS NOT found in the routing table of `root`.
295: The route 1.1.0.0/24 IS NOT found in the routing table of `root`.
296: The route 0.0.0.0/0 IS NOT found in the routing table of `net`.
297: The route 0.0.0.0/0 IS NOT found in the routing table of `pc`.
298: The route 3.1.0.0/24 IS NOT found in the routing table of `pc`.
299: The route 0.0.0.0/0 IS NOT found in the routing table of `local`.
300: The route 3.2.0.0/24 IS NOT found in the routing table of `local`.
301: named not started in the startup file of `root`
302: named not started in the startup file of `root`
303: named not started in the startup file of `local`
304: named not started in the startup file of `net`
305: `resolv.conf` file not found for device `as1r1`
306: `resolv.conf` file not found for device `as1r2`
307: `resolv.conf` file not found for device `as2r1`
308: `resolv.conf` file not found for device `as2r2`
309: `resolv.conf` file not found for device `as3r1`
310: `resolv.conf` file not found for device `as3r2`
311: `resolv.conf` file not found for device `pc`
312: ping: pc.net: Temporary failure in name resolution
313: ping: pc.net: Temporary failure in name resolution
314: ping: pc.net: Temporary failure in name resolution
315: ping: pc.net: Temporary failure in name resolution
316: ping: pc.net: Temporary failure in name resolution
317: ping: pc.net: Temporary failure in name resolution

'@

$eLab1 = @'
1: Device root not in the network scenario.
2: Device net not in the network scenario.
3: Devices connected to collision domain A ['as1r1', 'pc'] are different from the one in the template ['root', 'as1r2'].
4: Devices connected to collision domain E ['as2r2', 'dnsnet'] are different from the one in the template ['as2r2', 'net'].
5: Devices connected to collision domain K ['as1r2', 'local'] are different from the one in the template ['as3r2', 'local'].
6: Devices connected to collision domain J ['as3r2', 'dnsroot'] are different from the one in the template ['as3r2', 'pc'].
7: The interface `eth0` of `as1r1` has the following IP addresses: ['1.0.0.1/24']`.
8: The interface `eth1` of `as1r1` has the following IP addresses: ['10.20.0.1/30']`.
9: The interface `eth0` of `as1r2` has the following IP addresses: ['1.1.0.2/24']`.
10: The interface `eth1` of `as1r2` has the following IP addresses: ['10.20.1.1/30']`.
11: The interface `eth2` of `as1r2` has the following IP addresses: ['1.2.0.1/24']`.
12: The interface `eth2` of `as3r2` has the following IP addresses: ['3.1.0.1/24']`.
13: Interface eth`3` not found on `as3r2`
14: The interface `eth0` of `local` has the following IP addresses: ['1.2.0.2/24']`.
15: The interface `eth0` of `pc` has the following IP addresses: ['1.0.0.2/24']`.
16: No answer from `2.1.0.2` to `as1r1`.
17: No answer from `3.1.0.2` to `as1r1`.
18: No answer from `3.2.0.1` to `as1r1`.
19: No answer from `3.2.0.2` to `as1r1`.
20: No answer from `1.0.0.2` to `as1r2`.
21: No answer from `2.1.0.2` to `as1r2`.
22: No answer from `3.1.0.2` to `as1r2`.
23: No answer from `3.2.0.1` to `as1r2`.
24: No answer from `3.2.0.2` to `as1r2`.
25: No answer from `1.0.0.2` to `as2r1`.
26: No answer from `2.1.0.2` to `as2r1`.
27: No answer from `3.1.0.2` to `as2r1`.
28: No answer from `3.2.0.1` to `as2r1`.
29: No answer from `3.2.0.2` to `as2r1`.
30: No answer from `1.0.0.2` to `as2r2`.
31: No answer from `2.1.0.2` to `as2r2`.
32: No answer from `3.1.0.2` to `as2r2`.
33: No answer from `3.2.0.1` to `as2r2`.
34: No answer from `3.2.0.2` to `as2r2`.
35: No answer from `1.0.0.2` to `as3r1`.
36: No answer from `2.1.0.2` to `as3r1`.
37: No answer from `3.1.0.2` to `as3r1`.
38: ping: connect: Network is unreachable
39: ping: connect: Network is unreachable
40: No answer from `1.0.0.2` to `as3r2`.
41: No answer from `2.1.0.2` to `as3r2`.
42: No answer from `3.1.0.2` to `as3r2`.
43: ping: connect: Network is unreachable
44: ping: connect: Network is unreachable
45: ping: connect: Network is unreachable
46: ping: connect: Network is unreachable
47: ping: connect: Network is unreachable
48: ping: connect: Network is unreachable
49: ping: connect: Network is unreachable
50: ping: connect: Network is unreachable
51: ping: connect: Network is unreachable
52: ping: connect: Network is unreachable
53: ping: connect: Network is unreachable
54: ping: connect: Network is unreachable
55: ping: connect: Network is unreachable
56: ping: connect: Network is unreachable
57: ping: connect: Network is unreachable
58: ping: connect: Network is unreachable
59: ping: connect: Network is unreachable
60: ping: connect: Network is unreachable
61: ping: connect: Network is unreachable
62: ping: connect: Network is unreachable
63: ping: connect: Network is unreachable
64: ping: connect: Network is unreachable
65: ping: connect: Network is unreachable
66: ping: connect: Network is unreachable
67: Device `root` is not running.
68: Device `root` is not running.
69: Device `root` is not running.
70: Device `root` is not running.
71: Device `root` is not running.
72: Device `root` is not running.
73: Device `root` is not running.
74: Device `root` is not running.
75: Device `root` is not running.
76: Device `root` is not running.
77: Device `root` is not running.
78: Device `root` is not running.
79: Device `root` is not running.
80: Device `root` is not running.
81: Device `root` is not running.
82: Device `root` is not running.
83: Device `root` is not running.
84: Device `root` is not running.
85: Device `root` is not running.
86: Device `root` is not running.
87: Device `root` is not running.
88: Device `root` is not running.
89: Device `net` is not running.
90: Device `net` is not running.
91: Device `net` is not running.
92: Device `net` is not running.
93: Device `net` is not running.
94: Device `net` is not running.
95: Device `net` is not running.
96: Device `net` is not running.
97: Device `net` is not running.
98: Device `net` is not running.
99: Device `net` is not running.
100: Device `net` is not running.
101: Device `net` is not running.
102: Device `net` is not running.
103: Device `net` is not running.
104: Device `net` is not running.
105: Device `net` is not running.
106: Device `net` is not running.
107: Device `net` is not running.
108: Device `net` is not running.
109: Device `net` is not running.
110: Device `net` is not running.
111: ping: connect: Network is unreachable
112: ping: connect: Network is unreachable
113: ping: connect: Network is unreachable
114: ping: connect: Network is unreachable
115: ping: connect: Network is unreachable
116: ping: connect: Network is unreachable
117: ping: connect: Network is unreachable
118: ping: connect: Network is unreachable
119: ping: connect: Network is unreachable
120: ping: connect: Network is unreachable
121: ping: connect: Network is unreachable
122: ping: connect: Network is unreachable
123: ping: connect: Network is unreachable
124: ping: connect: Network is unreachable
125: ping: connect: Network is unreachable
126: ping: connect: Network is unreachable
127: ping: connect: Network is unreachable
128: ping: connect: Network is unreachable
129: ping: connect: Network is unreachable
130: ping: connect: Network is unreachable
131: Device net not in the network scenario.
132: Device net not in the network scenario.
133: Device root not in the network scenario.
134: Device root not in the network scenario.
135: The peering between as1r1 and 1.0.0.2 is not up.
136: The peering between as1r2 and 1.0.0.1 is not up.
137: The route 3.2.0.0/24 IS NOT found in the routing table of `as3r1`.
138: The route 0.0.0.0/0 IS NOT found in the routing table of `root`.
139: The route 1.1.0.0/24 IS NOT found in the routing table of `root`.
140: The route 0.0.0.0/0 IS NOT found in the routing table of `net`.
141: The route 2.1.0.0/24 IS NOT found in the routing table of `net`.
142: The route 0.0.0.0/0 IS NOT found in the routing table of `pc`.
143: The route 3.1.0.0/24 IS NOT found in the routing table of `pc`.
144: The route 0.0.0.0/0 IS NOT found in the routing table of `local`.
145: The route 3.2.0.0/24 IS NOT found in the routing table of `local`.
146: Device `root` is not running.
147: Device `root` is not running.
148: named on local is running but answered with REFUSED when quering for .
149: Device `net` is not running.
150: `resolv.conf` file not found for device `as1r1`
151: `resolv.conf` file not found for device `as1r2`
152: `resolv.conf` file not found for device `as2r1`
153: `resolv.conf` file not found for device `as2r2`
154: `resolv.conf` file not found for device `as3r1`
155: `resolv.conf` file not found for device `as3r2`
156: The local name server for device `pc` has ip `3.2.0.2`
157: ping: pc.net: Temporary failure in name resolution
158: ping: pc.net: Temporary failure in name resolution
159: ping: pc.net: Temporary failure in name resolution
160: ping: pc.net: Temporary failure in name resolution
161: ping: pc.net: Temporary failure in name resolution
162: ping: pc.net: Temporary failure in name resolution

'@

$ws.Range("A2").Value = "lab2"
$ws.Range("B2").Value = 387
$ws.Range("C2").Value = 14
$ws.Range("D2").Value = 401
$ws.Range("E2").Value = $eLab2
$ws.Range("E2").WrapText = $true

$ws.Range("A3").Value = "lab1"
$ws.Range("B3").Value = 237
$ws.Range("C3").Value = 162
$ws.Range("D3").Value = 399
$ws.Range("E3").Value = $eLab1
$ws.Range("E3").WrapText = $true

$ws.Range("A4").Value = "lab4"
$ws.Range("B4").Value = 401
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 401
$ws.Range("E4").Value = "None"
$ws.Range("E4").Style = "Normal"

$ws.Range("A5").Value = "lab3"
$ws.Range("B5").Value = 84
$ws.Range("C5").Value = 317
$ws.Range("D5").Value = 401
$ws.Range("E5").Value = $eLab3
$ws.Range("E5").WrapText = $true

$ws.Rows("2:5").AutoFit()
